$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrease Cheese, Tomatoe Sauce, and Dough inventory counts by 1
$ws.Range("B2").Value = 993
$ws.Range("C2").Value = 993
$ws.Range("D2").Value = 993

# Increase Sausage inventory count by 1
$ws.Range("G2").Value = 996
